$wb = $excel.ActiveWorkbook

# ----- About sheet -----
$about = $wb.Worksheets.Item("About")

# Remove the redundant/no-op style that used to sit on A8 (font-only xf identical to default)
$about.Range("A8").ClearFormats()

# Update the "desired output units" text block for the EU model
$about.Range("A11").Value = "For the EU. model, the desired output units are:"
$about.Range("A12").Value = "trillion passenger-km"
$about.Range("A13").Value = "trillion freight ton-km"

# New "Relevant Conversion Factors" section
$about.Range("A15").Value = "Relevant Conversion Factors"
$about.Range("A15:B15").Interior.Color = 12566463
$about.Range("A15").Font.Bold = $true

$about.Range("A16").Value = "miles to km"
$about.Range("B16").Value = 0.62137

# Column widths
$about.Columns.Item(1).ColumnWidth = 12.1
$about.Columns.Item(2).ColumnWidth = 15.6

# ----- CDCF-PMpPDOU sheet -----
$s2 = $wb.Worksheets.Item("CDCF-PMpPDOU")
$s2.Range("B2").Formula = "=10^12*About!B16"
$s2.Columns.Item(2).ColumnWidth = 27.6

# ----- CDCF-FTMpFDOU sheet -----
$s3 = $wb.Worksheets.Item("CDCF-FTMpFDOU")
$s3.Range("B2").Formula = "=10^12*About!B16"

# ----- Restore sheet selections (cosmetic, matches saved view state) -----
$about.Activate()
$about.Range("B17").Select()

$s2.Activate()
$s2.Range("B2").Select()

$s3.Activate()
$s3.Range("H13").Select()

$about.Activate()

